{"js": "// Replace the date line and each three-digit \u00f7 one-digit equation cell\n// in the practice-sheet table with the new values from the commit.\nconst replacements = [\n  [\"2025-10-11 Saturday\", \"2025-10-12 Sunday\"],\n  [\"370\u00f73=123, 1\", \"168\u00f76=28, 0\"],\n  [\"106\u00f74=26, 2\", \"227\u00f77=32, 3\"],\n  [\"750\u00f78=93, 6\", \"161\u00f72=80, 1\"],\n  [\"513\u00f75=102, 3\", \"263\u00f73=87, 2\"],\n  [\"560\u00f78=70, 0\", \"318\u00f77=45, 3\"],\n  [\"499\u00f72=249, 1\", \"377\u00f73=125, 2\"],\n  [\"994\u00f72=497, 0\", \"687\u00f74=171, 3\"],\n  [\"777\u00f76=129, 3\", \"785\u00f79=87, 2\"],\n  [\"809\u00f77=115, 4\", \"793\u00f79=88, 1\"],\n  [\"962\u00f76=160, 2\", \"485\u00f72=242, 1\"],\n  [\"243\u00f76=40, 3\", \"756\u00f76=126, 0\"],\n  [\"631\u00f74=157, 3\", \"383\u00f74=95, 3\"],\n  [\"969\u00f78=121, 1\", \"428\u00f77=61, 1\"],\n  [\"413\u00f72=206, 1\", \"114\u00f72=57, 0\"],\n  [\"260\u00f74=65, 0\", \"596\u00f73=198, 2\"],\n  [\"773\u00f78=96, 5\", \"157\u00f77=22, 3\"],\n  [\"382\u00f72=191, 0\", \"482\u00f73=160, 2\"],\n  [\"283\u00f72=141, 1\", \"405\u00f75=81, 0\"],\n  [\"949\u00f72=474, 1\", \"453\u00f72=226, 1\"],\n  [\"743\u00f79=82, 5\", \"870\u00f77=124, 2\"],\n  [\"743\u00f76=123, 5\", \"819\u00f75=163, 4\"],\n  [\"550\u00f74=137, 2\", \"327\u00f74=81, 3\"],\n  [\"983\u00f73=327, 2\", \"570\u00f76=95, 0\"],\n  [\"753\u00f78=94, 1\", \"332\u00f77=47, 3\"],\n  [\"440\u00f78=55, 0\", \"405\u00f74=101, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each three-digit \u00f7 one-digit equation cell\n# in the practice-sheet table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-11 Saturday\", \"2025-10-12 Sunday\"),\n    @(\"370\u00f73=123, 1\", \"168\u00f76=28, 0\"),\n    @(\"106\u00f74=26, 2\", \"227\u00f77=32, 3\"),\n    @(\"750\u00f78=93, 6\", \"161\u00f72=80, 1\"),\n    @(\"513\u00f75=102, 3\", \"263\u00f73=87, 2\"),\n    @(\"560\u00f78=70, 0\", \"318\u00f77=45, 3\"),\n    @(\"499\u00f72=249, 1\", \"377\u00f73=125, 2\"),\n    @(\"994\u00f72=497, 0\", \"687\u00f74=171, 3\"),\n    @(\"777\u00f76=129, 3\", \"785\u00f79=87, 2\"),\n    @(\"809\u00f77=115, 4\", \"793\u00f79=88, 1\"),\n    @(\"962\u00f76=160, 2\", \"485\u00f72=242, 1\"),\n    @(\"243\u00f76=40, 3\", \"756\u00f76=126, 0\"),\n    @(\"631\u00f74=157, 3\", \"383\u00f74=95, 3\"),\n    @(\"969\u00f78=121, 1\", \"428\u00f77=61, 1\"),\n    @(\"413\u00f72=206, 1\", \"114\u00f72=57, 0\"),\n    @(\"260\u00f74=65, 0\", \"596\u00f73=198, 2\"),\n    @(\"773\u00f78=96, 5\", \"157\u00f77=22, 3\"),\n    @(\"382\u00f72=191, 0\", \"482\u00f73=160, 2\"),\n    @(\"283\u00f72=141, 1\", \"405\u00f75=81, 0\"),\n    @(\"949\u00f72=474, 1\", \"453\u00f72=226, 1\"),\n    @(\"743\u00f79=82, 5\", \"870\u00f77=124, 2\"),\n    @(\"743\u00f76=123, 5\", \"819\u00f75=163, 4\"),\n    @(\"550\u00f74=137, 2\", \"327\u00f74=81, 3\"),\n    @(\"983\u00f73=327, 2\", \"570\u00f76=95, 0\"),\n    @(\"753\u00f78=94, 1\", \"332\u00f77=47, 3\"),\n    @(\"440\u00f78=55, 0\", \"405\u00f74=101, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
